$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1085.091
$ws.Range("I33").Value = 487.8
$ws.Range("J33").Value = 1582.8334
$ws.Range("K33").Value = 487.8
$ws.Range("L33").Value = 1582.8334
$ws.Range("M33").Value = -258.8
$ws.Range("N33").Value = -2040.8334
$ws.Range("H40").Value = 7464.0835
$ws.Range("I40").Value = 4906
$ws.Range("K40").Value = 4906
$ws.Range("M40").Value = -4731
$ws.Range("H53").Value = 437.92856
$ws.Range("I53").Value = 115.833336
$ws.Range("K53").Value = 115.833336
$ws.Range("M53").Value = 521.166664
$ws.Range("H62").Value = 2501.25
$ws.Range("I62").Value = 2501.25
$ws.Range("K62").Value = 2501.25
$ws.Range("M62").Value = -1877.25
$ws.Range("H64").Value = 7804
$ws.Range("I64").Value = 5950.857
$ws.Range("J64").Value = 9657.143
$ws.Range("K64").Value = 5950.857
$ws.Range("L64").Value = 9657.143
$ws.Range("M64").Value = -5702.857
$ws.Range("N64").Value = -10153.143
$ws.Range("H65").Value = 2501.25
$ws.Range("I65").Value = 2501.25
$ws.Range("K65").Value = 12506.25
$ws.Range("M65").Value = -9386.25
$ws.Range("H67").Value = 7804
$ws.Range("I67").Value = 5950.857
$ws.Range("J67").Value = 9657.143
$ws.Range("K67").Value = 5950.857
$ws.Range("L67").Value = 9657.143
$ws.Range("M67").Value = -5092.857
$ws.Range("N67").Value = -11373.143
$ws.Range("H82").Value = 2938.6667
$ws.Range("I82").Value = 2938.6667
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 8816.000100000001
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -8410.000100000001
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2938.6667
$ws.Range("I85").Value = 2938.6667
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 8816.000100000001
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -7412.000100000001
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 3603.0715
$ws.Range("I86").Value = 3086.125
$ws.Range("K86").Value = 3086.125
$ws.Range("M86").Value = -1963.125
$ws.Range("H88").Value = 6312.5835
$ws.Range("J88").Value = 5281.375
$ws.Range("L88").Value = 5281.375
$ws.Range("N88").Value = -6093.375
$ws.Range("H89").Value = 3603.0715
$ws.Range("I89").Value = 3086.125
$ws.Range("K89").Value = 15430.625
$ws.Range("M89").Value = -9814.625
$ws.Range("H91").Value = 6312.5835
$ws.Range("J91").Value = 5281.375
$ws.Range("L91").Value = 5281.375
$ws.Range("N91").Value = -8089.375
$ws.Range("H107").Value = 1611.4286
$ws.Range("I107").Value = 1974.5714
$ws.Range("K107").Value = 1974.5714
$ws.Range("M107").Value = -54.57140000000004
$ws.Range("H113").Value = 8141.357
$ws.Range("I113").Value = 5711.2856
$ws.Range("J113").Value = 10571.429
$ws.Range("K113").Value = 5711.2856
$ws.Range("L113").Value = 10571.429
$ws.Range("M113").Value = -2457.2856
$ws.Range("N113").Value = -17079.429
$ws.Range("H116").Value = 6400.8
$ws.Range("I116").Value = 5751
$ws.Range("K116").Value = 5751
$ws.Range("M116").Value = -2309
$ws.Range("H127").Value = 17772.7
$ws.Range("I127").Value = 1878.4
$ws.Range("J127").Value = 33667
$ws.Range("K127").Value = 5635.200000000001
$ws.Range("L127").Value = 101001
$ws.Range("M127").Value = -675.2000000000007
$ws.Range("N127").Value = -110921
$ws.Range("H129").Value = 113204.164
$ws.Range("I129").Value = 183288.73
$ws.Range("J129").Value = 3071.2856
$ws.Range("K129").Value = 549866.1900000001
$ws.Range("L129").Value = 9213.856800000001
$ws.Range("M129").Value = -544866.1900000001
$ws.Range("N129").Value = -19213.8568
$ws.Range("H132").Value = 12755.25
$ws.Range("I132").Value = 1975.3903
$ws.Range("K132").Value = 5926.1709
$ws.Range("M132").Value = -3396.1709
$ws.Range("H137").Value = 3788.0476
$ws.Range("I137").Value = 3908.8333
$ws.Range("J137").Value = 3063.3333
$ws.Range("K137").Value = 11726.4999
$ws.Range("L137").Value = 9189.999899999999
$ws.Range("M137").Value = -9176.499899999999
$ws.Range("N137").Value = -14289.9999
$ws.Range("H141").Value = 4967.8965
$ws.Range("I141").Value = 2079.577
$ws.Range("J141").Value = 30000
$ws.Range("K141").Value = 6238.731000000001
$ws.Range("L141").Value = 90000
$ws.Range("M141").Value = -1058.731000000001
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15388698
$ws.Range("I32").Value = 15388698
$ws.Range("K32").Value = 15388698
$ws.Range("M32").Value = -15388411
$ws.Range("H45").Value = 5840.6665
$ws.Range("I45").Value = 2348.3333
$ws.Range("J45").Value = 9333
$ws.Range("K45").Value = 2348.3333
$ws.Range("L45").Value = 9333
$ws.Range("M45").Value = -1971.3333
$ws.Range("N45").Value = -10087
$ws.Range("H61").Value = 2149.2856
$ws.Range("I61").Value = 2176.75
$ws.Range("K61").Value = 2176.75
$ws.Range("M61").Value = -1964.75
$ws.Range("H63").Value = 11149.167
$ws.Range("I63").Value = 6723.75
$ws.Range("K63").Value = 6723.75
$ws.Range("M63").Value = -6037.75
$ws.Range("H66").Value = 11149.167
$ws.Range("I66").Value = 6723.75
$ws.Range("K66").Value = 33618.75
$ws.Range("M66").Value = -30186.75
$ws.Range("H74").Value = 2073.875
$ws.Range("I74").Value = 2073.875
$ws.Range("K74").Value = 2073.875
$ws.Range("M74").Value = -1199.875
$ws.Range("H77").Value = 2073.875
$ws.Range("I77").Value = 2073.875
$ws.Range("K77").Value = 10369.375
$ws.Range("M77").Value = -6001.375
$ws.Range("H88").Value = 1204.2222
$ws.Range("I88").Value = 893.4
$ws.Range("J88").Value = 1323.7693
$ws.Range("K88").Value = 893.4
$ws.Range("L88").Value = 1323.7693
$ws.Range("M88").Value = -487.4
$ws.Range("N88").Value = -2135.7693
$ws.Range("H91").Value = 1204.2222
$ws.Range("I91").Value = 893.4
$ws.Range("J91").Value = 1323.7693
$ws.Range("K91").Value = 893.4
$ws.Range("L91").Value = 1323.7693
$ws.Range("M91").Value = 510.6
$ws.Range("N91").Value = -4131.7693
$ws.Range("H102").Value = 1415
$ws.Range("I102").Value = 1415
$ws.Range("K102").Value = 1415
$ws.Range("M102").Value = 207
$ws.Range("H110").Value = 7960.5
$ws.Range("I110").Value = 7975.6
$ws.Range("K110").Value = 7975.6
$ws.Range("M110").Value = -5930.6
$ws.Range("H122").Value = 3433.75
$ws.Range("I122").Value = 3433.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10301.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7851.25
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1683.2894
$ws.Range("I132").Value = 1674.7297
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5024.189100000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2494.189100000001
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 2149.2856
$ws.Range("I136").Value = 2176.75
$ws.Range("K136").Value = 6530.25
$ws.Range("M136").Value = -3980.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1424.8572
$ws.Range("I20").Value = 1433.3334
$ws.Range("J20").Value = 1418.5
$ws.Range("K20").Value = 1433.3334
$ws.Range("L20").Value = 1418.5
$ws.Range("M20").Value = -1186.3334
$ws.Range("N20").Value = -1912.5
$ws.Range("H22").Value = 269.33334
$ws.Range("I22").Value = 269.33334
$ws.Range("K22").Value = 269.33334
$ws.Range("M22").Value = -96.33334000000002
$ws.Range("H63").Value = 50271
$ws.Range("J63").Value = 50271
$ws.Range("L63").Value = 50271
$ws.Range("N63").Value = -51643
$ws.Range("H66").Value = 50271
$ws.Range("J66").Value = 50271
$ws.Range("L66").Value = 150813
$ws.Range("N66").Value = -157677
$ws.Range("H107").Value = 19995
$ws.Range("I107").Value = 9990
$ws.Range("J107").Value = 30000
$ws.Range("K107").Value = 9990
$ws.Range("L107").Value = 30000
$ws.Range("M107").Value = -8070
$ws.Range("N107").Value = -33840
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 2366.9048
$ws.Range("I134").Value = 2135.25
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 6405.75
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -3870.75
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 17499
$ws.Range("I16").Value = 18873.75
$ws.Range("J16").Value = 12000
$ws.Range("K16").Value = 18873.75
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = -18586.75
$ws.Range("N16").Value = -12574
$ws.Range("H31").Value = 1831.5143
$ws.Range("I31").Value = 1798.5714
$ws.Range("K31").Value = 1798.5714
$ws.Range("M31").Value = -1503.5714
$ws.Range("H33").Value = 3004.1428
$ws.Range("I33").Value = 2619
$ws.Range("J33").Value = 3967
$ws.Range("K33").Value = 2619
$ws.Range("L33").Value = 3967
$ws.Range("M33").Value = -2240
$ws.Range("N33").Value = -4725
$ws.Range("H34").Value = 1831.5143
$ws.Range("I34").Value = 1798.5714
$ws.Range("K34").Value = 1798.5714
$ws.Range("M34").Value = -1596.5714
$ws.Range("H58").Value = 1336.421
$ws.Range("I58").Value = 753.25
$ws.Range("J58").Value = 4446.6665
$ws.Range("K58").Value = 753.25
$ws.Range("L58").Value = 4446.6665
$ws.Range("M58").Value = -550.25
$ws.Range("N58").Value = -4852.6665
$ws.Range("H99").Value = 32063952
$ws.Range("I99").Value = 8133881.5
$ws.Range("K99").Value = 8133881.5
$ws.Range("M99").Value = -8132383.5
$ws.Range("H105").Value = 3171.2666
$ws.Range("J105").Value = 2170
$ws.Range("L105").Value = 2170
$ws.Range("N105").Value = -5664
$ws.Range("H107").Value = 10519.81
$ws.Range("J107").Value = 14395.066
$ws.Range("L107").Value = 14395.066
$ws.Range("N107").Value = -18235.066
$ws.Range("H113").Value = 17499
$ws.Range("I113").Value = 18873.75
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 18873.75
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -16703.75
$ws.Range("N113").Value = -16340
$ws.Range("H126").Value = 32063952
$ws.Range("I126").Value = 8133881.5
$ws.Range("K126").Value = 24401644.5
$ws.Range("M126").Value = -24399174.5
$ws.Range("H132").Value = 1717.037
$ws.Range("I132").Value = 1474.4
$ws.Range("K132").Value = 4423.200000000001
$ws.Range("M132").Value = -1893.200000000001
$ws.Range("H134").Value = 3104.5925
$ws.Range("I134").Value = 2616.4
$ws.Range("J134").Value = 9207
$ws.Range("K134").Value = 7849.200000000001
$ws.Range("L134").Value = 27621
$ws.Range("M134").Value = -5314.200000000001
$ws.Range("N134").Value = -32691
$ws.Range("H136").Value = 1336.421
$ws.Range("I136").Value = 753.25
$ws.Range("J136").Value = 4446.6665
$ws.Range("K136").Value = 2259.75
$ws.Range("L136").Value = 13339.9995
$ws.Range("M136").Value = 290.25
$ws.Range("N136").Value = -18439.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1719.8
$ws.Range("J7").Value = 3975
$ws.Range("L7").Value = 11925
$ws.Range("N7").Value = -12149
$ws.Range("H126").Value = 2249.5
$ws.Range("H131").Value = 3422.3684
$ws.Range("I131").Value = 2029.6
$ws.Range("K131").Value = 6088.799999999999
$ws.Range("M131").Value = -1048.799999999999
$ws.Range("H137").Value = 2863.68
$ws.Range("J137").Value = 3115.3572
$ws.Range("L137").Value = 9346.071599999999
$ws.Range("N137").Value = -19546.0716
$ws.Range("H138").Value = 5386.615
$ws.Range("I138").Value = 3167.077
$ws.Range("J138").Value = 7606.154
$ws.Range("K138").Value = 9501.231
$ws.Range("L138").Value = 22818.462
$ws.Range("M138").Value = -4361.231
$ws.Range("N138").Value = -33098.462
$ws.Range("H139").Value = 2200
$ws.Range("J139").Value = 3000
$ws.Range("L139").Value = 9000
$ws.Range("N139").Value = -19280
$ws.Range("H140").Value = 4219
$ws.Range("I140").Value = 2026.5
$ws.Range("J140").Value = 6850
$ws.Range("K140").Value = 6079.5
$ws.Range("L140").Value = 20550
$ws.Range("M140").Value = -899.5
$ws.Range("N140").Value = -30910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8026.423
$ws.Range("I80").Value = 7051.5625
$ws.Range("K80").Value = 7051.5625
$ws.Range("M80").Value = -6053.5625
$ws.Range("H83").Value = 8026.423
$ws.Range("I83").Value = 7051.5625
$ws.Range("K83").Value = 35257.8125
$ws.Range("M83").Value = -30265.8125
$ws.Range("H97").Value = 457.375
$ws.Range("I97").Value = 401.46155
$ws.Range("K97").Value = 401.46155
$ws.Range("M97").Value = 94.53845000000001
$ws.Range("H113").Value = 8850.799999999999
$ws.Range("I113").Value = 8127
$ws.Range("J113").Value = 9333.333000000001
$ws.Range("K113").Value = 8127
$ws.Range("L113").Value = 9333.333000000001
$ws.Range("M113").Value = -5957
$ws.Range("N113").Value = -13673.333
$ws.Range("H117").Value = 59892.5
$ws.Range("J117").Value = 59892.5
$ws.Range("L117").Value = 59892.5
$ws.Range("N117").Value = -66776.5
$ws.Range("H122").Value = 3690.8667
$ws.Range("I122").Value = 3165.1738
$ws.Range("J122").Value = 5418.143
$ws.Range("K122").Value = 9495.5214
$ws.Range("L122").Value = 16254.429
$ws.Range("M122").Value = -7045.5214
$ws.Range("N122").Value = -21154.429
$ws.Range("H132").Value = 5195
$ws.Range("I132").Value = 6045
$ws.Range("J132").Value = 3495
$ws.Range("K132").Value = 18135
$ws.Range("L132").Value = 10485
$ws.Range("M132").Value = -15605
$ws.Range("N132").Value = -15545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5574.3447
$ws.Range("I7").Value = 5166.7144
$ws.Range("J7").Value = 6644.375
$ws.Range("K7").Value = 5166.7144
$ws.Range("L7").Value = 6644.375
$ws.Range("M7").Value = -5054.7144
$ws.Range("N7").Value = -6868.375
$ws.Range("H12").Value = 263
$ws.Range("J12").Value = 299.5
$ws.Range("L12").Value = 299.5
$ws.Range("N12").Value = -639.5
$ws.Range("H16").Value = 786.625
$ws.Range("I16").Value = 799
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 799
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -629
$ws.Range("N16").Value = -1040
$ws.Range("H46").Value = 2961.6667
$ws.Range("I46").Value = 2900
$ws.Range("J46").Value = 2992.5
$ws.Range("K46").Value = 2900
$ws.Range("L46").Value = 2992.5
$ws.Range("M46").Value = -2712
$ws.Range("N46").Value = -3368.5
$ws.Range("H61").Value = 1098.08
$ws.Range("I61").Value = 803.2174
$ws.Range("J61").Value = 4489
$ws.Range("K61").Value = 803.2174
$ws.Range("L61").Value = 4489
$ws.Range("M61").Value = -601.2174
$ws.Range("N61").Value = -4893
$ws.Range("H68").Value = 3280.6924
$ws.Range("I68").Value = 3066
$ws.Range("J68").Value = 3996.3333
$ws.Range("K68").Value = 3066
$ws.Range("L68").Value = 3996.3333
$ws.Range("M68").Value = -2317
$ws.Range("N68").Value = -5494.3333
$ws.Range("H69").Value = 59000
$ws.Range("J69").Value = 59000
$ws.Range("L69").Value = 59000
$ws.Range("N69").Value = -60622
$ws.Range("H71").Value = 3280.6924
$ws.Range("I71").Value = 3066
$ws.Range("J71").Value = 3996.3333
$ws.Range("K71").Value = 15330
$ws.Range("L71").Value = 19981.6665
$ws.Range("M71").Value = -11586
$ws.Range("N71").Value = -27469.6665
$ws.Range("H72").Value = 59000
$ws.Range("J72").Value = 59000
$ws.Range("L72").Value = 177000
$ws.Range("N72").Value = -185112
$ws.Range("H113").Value = 1098.08
$ws.Range("I113").Value = 803.2174
$ws.Range("J113").Value = 4489
$ws.Range("K113").Value = 803.2174
$ws.Range("L113").Value = 4489
$ws.Range("M113").Value = 1366.7826
$ws.Range("N113").Value = -8829
$ws.Range("H122").Value = 5722.68
$ws.Range("I122").Value = 2770.25
$ws.Range("K122").Value = 8310.75
$ws.Range("M122").Value = -5860.75
$ws.Range("H123").Value = 57500
$ws.Range("J123").Value = 57500
$ws.Range("L123").Value = 57500
$ws.Range("N123").Value = -67300
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840
$ws.Range("H126").Value = 5574.3447
$ws.Range("I126").Value = 5166.7144
$ws.Range("J126").Value = 6644.375
$ws.Range("K126").Value = 15500.1432
$ws.Range("L126").Value = 19933.125
$ws.Range("M126").Value = -13030.1432
$ws.Range("N126").Value = -24873.125
$ws.Range("H132").Value = 8228.429
$ws.Range("I132").Value = 2519.8
$ws.Range("J132").Value = 22500
$ws.Range("K132").Value = 7559.400000000001
$ws.Range("L132").Value = 67500
$ws.Range("M132").Value = -5029.400000000001
$ws.Range("N132").Value = -72560
$ws.Range("H136").Value = 3076.5
$ws.Range("I136").Value = 3125.68
$ws.Range("J136").Value = 2666.6667
$ws.Range("K136").Value = 9377.039999999999
$ws.Range("L136").Value = 8000.000100000001
$ws.Range("M136").Value = -6827.039999999999
$ws.Range("N136").Value = -13100.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20016990
$ws.Range("I2").Value = 33353334
$ws.Range("K2").Value = 33353334
$ws.Range("M2").Value = -33353222
$ws.Range("H62").Value = 2333
$ws.Range("I62").Value = 2499.5
$ws.Range("K62").Value = 2499.5
$ws.Range("M62").Value = -1875.5
$ws.Range("H65").Value = 2333
$ws.Range("I65").Value = 2499.5
$ws.Range("K65").Value = 12497.5
$ws.Range("M65").Value = -9377.5
$ws.Range("H100").Value = 1014.0714
$ws.Range("I100").Value = 1115.091
$ws.Range("J100").Value = 643.6667
$ws.Range("K100").Value = 2230.182
$ws.Range("L100").Value = 1287.3334
$ws.Range("M100").Value = -1689.182
$ws.Range("N100").Value = -2369.3334
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 3589.8
$ws.Range("I122").Value = 3589.8
$ws.Range("K122").Value = 10769.4
$ws.Range("M122").Value = -8319.400000000001
$ws.Range("H126").Value = 1590.8462
$ws.Range("I126").Value = 1463.1428
$ws.Range("J126").Value = 2127.2
$ws.Range("K126").Value = 4389.428400000001
$ws.Range("L126").Value = 6381.599999999999
$ws.Range("M126").Value = -1919.428400000001
$ws.Range("N126").Value = -11321.6
$ws.Range("H132").Value = 1516.4736
$ws.Range("I132").Value = 1459.3889
$ws.Range("K132").Value = 4378.1667
$ws.Range("M132").Value = -1848.1667
$ws.Range("H136").Value = 3012.8
$ws.Range("I136").Value = 1836.75
$ws.Range("K136").Value = 5510.25
$ws.Range("M136").Value = -2960.25
